# Horarios actualizados Linea 141 - 310
# Applies the scrape-refresh update to all three sheets:
#   - "Ultima actualizacion" timestamp -> 04:54:25
#   - "Total filas" counts bumped
#   - Sheet "LP1912": a new arrival (11_ETCHEVERRY) is inserted at row 15
#     (pushing the previous rows 15-28 down to 16-29), and two brand new
#     arrivals are appended at the end (rows 30-31).
#   - Sheet "LP1912-215": one new arrival appended at the end (row 15).
#   - Sheet "6203-6173": one new arrival appended at the end (row 12).

$wb = $excel.ActiveWorkbook

$newUpdateTime = "04:54:25"

function Set-Header($ws, [string]$totalFilasText) {
    $ws.Range("A2").Value = "Última actualización: $newUpdateTime"
    $ws.Range("A3").Value = $totalFilasText
}

function Set-DataRow($ws, [int]$row, [string]$horaScrap, [string]$horaLlegada, [string]$linea, $minutos, [string]$parada) {
    $ws.Cells.Item($row, 1).Value = $horaScrap
    $ws.Cells.Item($row, 2).Value = $horaLlegada
    $ws.Cells.Item($row, 3).Value = $linea
    $ws.Cells.Item($row, 4).Value = $minutos
    $ws.Cells.Item($row, 5).Value = $parada
}

# ---------------------------------------------------------------------------
# Sheet 1: LP1912  (23 -> 26 rows)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

Set-Header $ws1 "Total filas: 26"

# Insert the new row at position 15, pushing existing data rows 15-28 down to 16-29
$ws1.Rows.Item(15).Insert()
Set-DataRow $ws1 15 "04:54:25" "04:54" "11_ETCHEVERRY" 0 "LP1912"

# Append two brand-new rows at the end
Set-DataRow $ws1 30 "04:54:25" "06:44" "225_C ROCA-H SUR" 110 "LP1912"
Set-DataRow $ws1 31 "04:54:25" "06:46" "215C_EL PATO" 112 "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215  (9 -> 10 rows)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

Set-Header $ws2 "Total filas: 10"

Set-DataRow $ws2 15 "04:54:25" "06:46" "215C_EL PATO" 112 "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173  (6 -> 7 rows)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

Set-Header $ws3 "Total filas: 7"

Set-DataRow $ws3 12 "04:54:25" "06:33" "215C_LA PLATA" 99 "L6203"

Write-Output "Horarios actualizados Linea 141 - 310 aplicado"
